$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.042552731154559
$ws.Range("D2").Value = 1.048823031399251
$ws.Range("E2").Value = 1.050467562864365
$ws.Range("F2").Value = 1.060805173670996
$ws.Range("I2").Value = 1.041624720762456
$ws.Range("J2").Value = 1.047627936952472
$ws.Range("K2").Value = 1.051581803157005
$ws.Range("L2").Value = 1.053221758389633
$ws.Range("M2").Value = 1.063530954634238
$ws.Range("N2").Value = 1.049115688787549
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.043386590563229
$ws.Range("D3").Value = 1.049466880285963
$ws.Range("E3").Value = 1.051205231418866
$ws.Range("F3").Value = 1.061613989357125
$ws.Range("I3").Value = 1.041805540576699
$ws.Range("J3").Value = 1.04810885633163
$ws.Range("K3").Value = 1.052038209965285
$ws.Range("L3").Value = 1.053772071986228
$ws.Range("M3").Value = 1.064154281893712
$ws.Range("N3").Value = 1.049597291127387
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.043926889628584
$ws.Range("D4").Value = 1.049884145871266
$ws.Range("E4").Value = 1.051683594754179
$ws.Range("F4").Value = 1.062138514721563
$ws.Range("I4").Value = 1.041921714711219
$ws.Range("J4").Value = 1.048420099927032
$ws.Range("K4").Value = 1.052333498719141
$ws.Range("L4").Value = 1.054128533716248
$ws.Range("M4").Value = 1.06455812664023
$ws.Range("N4").Value = 1.049908976724395
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.044154205852373
$ws.Range("D5").Value = 1.050059718623298
$ws.Range("E5").Value = 1.051884946317846
$ws.Range("F5").Value = 1.062359302622904
$ws.Range("I5").Value = 1.041970355263864
$ws.Range("J5").Value = 1.048550958932383
$ws.Range("K5").Value = 1.052457627888393
$ws.Range("L5").Value = 1.054278477770087
$ws.Range("M5").Value = 1.064728023437885
$ws.Range("N5").Value = 1.050040021564543
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.044192383427835
$ws.Range("D6").Value = 1.050089207054955
$ws.Range("E6").Value = 1.051918768595119
$ws.Range("F6").Value = 1.062396390091922
$ws.Range("I6").Value = 1.041978510532655
$ws.Range("J6").Value = 1.048572931411175
$ws.Range("K6").Value = 1.05247846908274
$ws.Range("L6").Value = 1.054303659125754
$ws.Range("M6").Value = 1.064756556872338
$ws.Range("N6").Value = 1.050062025246776
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.043929926355648
$ws.Range("D7").Value = 1.049886491279021
$ws.Range("E7").Value = 1.051686284253109
$ws.Range("F7").Value = 1.062141463812766
$ws.Range("I7").Value = 1.041922365431727
$ws.Range("J7").Value = 1.048421848425232
$ws.Range("K7").Value = 1.052335157381066
$ws.Range("L7").Value = 1.054130536933964
$ws.Range("M7").Value = 1.06456039633835
$ws.Range("N7").Value = 1.049910727705663
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.042834384332627
$ws.Range("D8").Value = 1.049040486962045
$ws.Range("E8").Value = 1.050716644388157
$ws.Range("F8").Value = 1.061078273728679
$ws.Range("I8").Value = 1.041686000741511
$ws.Range("J8").Value = 1.047790453376092
$ws.Range("K8").Value = 1.051736054804636
$ws.Range("L8").Value = 1.05340766149813
$ws.Range("M8").Value = 1.063741504306405
$ws.Range("N8").Value = 1.049278436003131
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.040909616565048
$ws.Range("D9").Value = 1.047554797481317
$ws.Range("E9").Value = 1.049016085175072
$ws.Range("F9").Value = 1.059213830494937
$ws.Range("I9").Value = 1.041263184323331
$ws.Range("J9").Value = 1.046678347938671
$ws.Range("K9").Value = 1.050680137085036
$ws.Range("L9").Value = 1.052136779351701
$ws.Range("M9").Value = 1.062302495977007
$ws.Range("N9").Value = 1.048164751248428
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.039630387623587
$ws.Range("D10").Value = 1.046567863481803
$ws.Range("E10").Value = 1.047887921410671
$ws.Range("F10").Value = 1.057977065864419
$ws.Range("I10").Value = 1.040977107371949
$ws.Range("J10").Value = 1.045937356798702
$ws.Range("K10").Value = 1.049976124424662
$ws.Range("L10").Value = 1.051291576696712
$ws.Range("M10").Value = 1.061345937854756
$ws.Range("N10").Value = 1.047422707816048
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.039077426438807
$ws.Range("D11").Value = 1.046141369995753
$ws.Range("E11").Value = 1.047400753306565
$ws.Range("F11").Value = 1.057443029345399
$ws.Range("I11").Value = 1.040852246217519
$ws.Range("J11").Value = 1.04561661437853
$ws.Range("K11").Value = 1.049671280515786
$ws.Range("L11").Value = 1.050926100624227
$ws.Range("M11").Value = 1.060932418808136
$ws.Range("N11").Value = 1.04710150990485
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.038872177066941
$ws.Range("D12").Value = 1.04598308176618
$ws.Range("E12").Value = 1.047219999848713
$ws.Range("F12").Value = 1.057244890507668
$ws.Range("I12").Value = 1.04080571954087
$ws.Range("J12").Value = 1.045497494513787
$ws.Range("K12").Value = 1.049558049049363
$ws.Range("L12").Value = 1.050790423558678
$ws.Range("M12").Value = 1.060778922878274
$ws.Range("N12").Value = 1.046982220876235
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.038916197162017
$ws.Range("D13").Value = 1.046017029195234
$ws.Range("E13").Value = 1.047258762884667
$ws.Range("F13").Value = 1.057287381694981
$ws.Range("I13").Value = 1.040815706342028
$ws.Range("J13").Value = 1.045523045285605
$ws.Range("K13").Value = 1.049582337509433
$ws.Range("L13").Value = 1.050819523239379
$ws.Range("M13").Value = 1.060811843566518
$ws.Range("N13").Value = 1.04700780793308
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.039060457492671
$ws.Range("D14").Value = 1.046128283169275
$ws.Range("E14").Value = 1.047385808027461
$ws.Range("F14").Value = 1.057426646502459
$ws.Range("I14").Value = 1.040848403316325
$ws.Range("J14").Value = 1.045606767515094
$ws.Range("K14").Value = 1.049661920736196
$ws.Range("L14").Value = 1.050914883935433
$ws.Range("M14").Value = 1.060919728666075
$ws.Range("N14").Value = 1.047091649057737
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.039149360315635
$ws.Range("D15").Value = 1.046196847756561
$ws.Range("E15").Value = 1.047464111621941
$ws.Range("F15").Value = 1.057512482171503
$ws.Range("I15").Value = 1.040868529453972
$ws.Range("J15").Value = 1.045658353990337
$ws.Range("K15").Value = 1.049710954781053
$ws.Range("L15").Value = 1.050973649060275
$ws.Range("M15").Value = 1.060986213997961
$ws.Range("N15").Value = 1.047143308791693
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.03966710607639
$ws.Range("D16").Value = 1.04659618663705
$ws.Range("E16").Value = 1.047920281436126
$ws.Range("F16").Value = 1.058012539699858
$ws.Range("I16").Value = 1.040985373244057
$ws.Range("J16").Value = 1.045958645886522
$ws.Range("K16").Value = 1.049996355993072
$ws.Range("L16").Value = 1.051315842855292
$ws.Range("M16").Value = 1.061373396163958
$ws.Range("N16").Value = 1.047444027136815
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.039992130830015
$ws.Range("D17").Value = 1.046846911749844
$ws.Range("E17").Value = 1.048206783345737
$ws.Range("F17").Value = 1.058326613196068
$ws.Range("I17").Value = 1.041058402289376
$ws.Range("J17").Value = 1.046147041920525
$ws.Range("K17").Value = 1.050175381146101
$ws.Range("L17").Value = 1.051530627553268
$ws.Range("M17").Value = 1.061616447780289
$ws.Range("N17").Value = 1.047632690714795
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.040181804137289
$ws.Range("D18").Value = 1.046993237875987
$ws.Range("E18").Value = 1.048374023663444
$ws.Range("F18").Value = 1.058509950642927
$ws.Range("I18").Value = 1.041100903523471
$ws.Range("J18").Value = 1.046256940818203
$ws.Range("K18").Value = 1.050279803223109
$ws.Range("L18").Value = 1.051655956189251
$ws.Range("M18").Value = 1.061758280872533
$ws.Range("N18").Value = 1.047742745681513
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.040246493361172
$ws.Range("D19").Value = 1.0470431452233
$ws.Range("E19").Value = 1.048431070041666
$ws.Range("F19").Value = 1.058572488281721
$ws.Range("I19").Value = 1.041115379146552
$ws.Range("J19").Value = 1.046294415271555
$ws.Range("K19").Value = 1.050315408365105
$ws.Range("L19").Value = 1.051698698169732
$ws.Range("M19").Value = 1.061806653296773
$ws.Range("N19").Value = 1.047780273352887
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.039957249234713
$ws.Range("D20").Value = 1.046820002780209
$ws.Range("E20").Value = 1.04817603108595
$ws.Range("F20").Value = 1.058292901204313
$ws.Range("I20").Value = 1.041050576823422
$ws.Range("J20").Value = 1.046126827694569
$ws.Range("K20").Value = 1.050156173458849
$ws.Range("L20").Value = 1.051507578174971
$ws.Range("M20").Value = 1.061590363885838
$ws.Range("N20").Value = 1.047612447782319
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.039017972402826
$ws.Range("D21").Value = 1.046095518030647
$ws.Range("E21").Value = 1.047348390794198
$ws.Range("F21").Value = 1.057385630229013
$ws.Range("I21").Value = 1.04083877894318
$ws.Range("J21").Value = 1.04558211289309
$ws.Range("K21").Value = 1.049638485416897
$ws.Range("L21").Value = 1.050886800448557
$ws.Range("M21").Value = 1.060887956323243
$ws.Range("N21").Value = 1.047066959423343
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.0384282517857
$ws.Range("D22").Value = 1.04564076135232
$ws.Range("E22").Value = 1.046829192900638
$ws.Range("F22").Value = 1.056816502179635
$ws.Range("I22").Value = 1.040704759122327
$ws.Range("J22").Value = 1.045239734777878
$ws.Range("K22").Value = 1.049313001655866
$ws.Range("L22").Value = 1.050496939243843
$ws.Range("M22").Value = 1.060446924242666
$ws.Range("N22").Value = 1.046724095091935
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.038740793618827
$ws.Range("D23").Value = 1.04588176420571
$ws.Range("E23").Value = 1.047104317725875
$ws.Range("F23").Value = 1.057118082913965
$ws.Range("I23").Value = 1.040775886273976
$ws.Range("J23").Value = 1.045421225414875
$ws.Range("K23").Value = 1.049485545599805
$ws.Range("L23").Value = 1.050703569210607
$ws.Range("M23").Value = 1.0606806662882
$ws.Range("N23").Value = 1.046905843466454
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.03997301044254
$ws.Range("D24").Value = 1.046832161530928
$ws.Range("E24").Value = 1.048189926311345
$ws.Range("F24").Value = 1.058308133758577
$ws.Range("I24").Value = 1.041054113109424
$ws.Range("J24").Value = 1.046135961600351
$ws.Range("K24").Value = 1.050164852586937
$ws.Range("L24").Value = 1.051517993047669
$ws.Range("M24").Value = 1.061602149873757
$ws.Range("N24").Value = 1.047621594659295
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.041406526776317
$ws.Range("D25").Value = 1.047938270548629
$ws.Range("E25").Value = 1.049454752135476
$ws.Range("F25").Value = 1.059694750472749
$ws.Range("I25").Value = 1.041373235911633
$ws.Range("J25").Value = 1.046965787066412
$ws.Range("K25").Value = 1.050953135032693
$ws.Range("L25").Value = 1.052464977817491
$ws.Range("M25").Value = 1.062674031469271
$ws.Range("N25").Value = 1.048452598572701
